# Swap the match-data (columns B through AD) between row 20 and row 21.
# Column A (the running index, 18/19) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Singapore Premier League")

$row1 = 20
$row2 = 21

# Columns B (2) .. AD (30)
$firstCol = 2
$lastCol = 30

$range1 = $ws.Range($ws.Cells.Item($row1, $firstCol), $ws.Cells.Item($row1, $lastCol))
$range2 = $ws.Range($ws.Cells.Item($row2, $firstCol), $ws.Cells.Item($row2, $lastCol))

$values1 = $range1.Value2
$values2 = $range2.Value2

$range1.Value2 = $values2
$range2.Value2 = $values1
